# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "K" column (column G) values per regenerated save_data
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("G8").Value = 4
$ws.Range("G9").Value = 2
$ws.Range("G10").Value = 1
$ws.Range("G12").Value = 1
$ws.Range("G13").Value = 1
